# Reverse the Mid AI heuristic matrix on the "Problem1" sheet.
# Every existing value of 0.5 becomes 0, and every existing value of 0 becomes -1
# (equivalent to: new = old*2 - 1), over the data range A2:J11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Problem1")

for ($r = 2; $r -le 11; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $old = $cell.Value2
        $cell.Value2 = ($old * 2) - 1
    }
}

# Update the stored cursor/selection position to match the author's edit.
$ws.Activate()
$ws.Range("J19").Select()
